{"js": "// Positional replacements: the document has one date paragraph and a table\n// of 20 rows x 5 columns where 5 rows (0, 4, 9, 14, 19) contain 25 total\n// \"NN\u00d7NN=\" multiplication prompts. Each old string below is replaced with\n// its corresponding new string, in document order. Because one string\n// (\"87\u00d761=\") occurs twice with two different replacements, we search again\n// after every replacement so the next hit (the first remaining unreplaced\n// occurrence) is the one that gets the next value in the list.\nconst replacements = [\n  [\"2024-08-22 Thursday\", \"2024-08-23 Friday\"],\n  [\"38\u00d769=\", \"50\u00d722=\"],\n  [\"65\u00d767=\", \"68\u00d793=\"],\n  [\"22\u00d798=\", \"38\u00d733=\"],\n  [\"45\u00d713=\", \"59\u00d728=\"],\n  [\"38\u00d730=\", \"52\u00d770=\"],\n  [\"48\u00d753=\", \"91\u00d742=\"],\n  [\"23\u00d766=\", \"79\u00d766=\"],\n  [\"80\u00d740=\", \"42\u00d798=\"],\n  [\"27\u00d770=\", \"35\u00d725=\"],\n  [\"35\u00d723=\", \"79\u00d766=\"],\n  [\"54\u00d753=\", \"56\u00d720=\"],\n  [\"87\u00d761=\", \"48\u00d743=\"],\n  [\"64\u00d798=\", \"29\u00d758=\"],\n  [\"87\u00d761=\", \"19\u00d727=\"],\n  [\"35\u00d791=\", \"42\u00d776=\"],\n  [\"72\u00d727=\", \"14\u00d788=\"],\n  [\"95\u00d792=\", \"98\u00d777=\"],\n  [\"12\u00d725=\", \"70\u00d776=\"],\n  [\"43\u00d711=\", \"84\u00d724=\"],\n  [\"20\u00d738=\", \"11\u00d794=\"],\n  [\"19\u00d769=\", \"51\u00d759=\"],\n  [\"23\u00d737=\", \"83\u00d751=\"],\n  [\"52\u00d793=\", \"29\u00d748=\"],\n  [\"42\u00d772=\", \"97\u00d746=\"],\n  [\"51\u00d742=\", \"61\u00d736=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  // The first hit is always the left-most not-yet-updated occurrence,\n  // because earlier occurrences of duplicated source text were already\n  // rewritten to their (different) target text in previous iterations.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Positional replacements: the document has one date paragraph and a table\n# of 20 rows x 5 columns where 5 rows (0, 4, 9, 14, 19) contain 25 total\n# \"NN\u00d7NN=\" multiplication prompts. Each old string below is replaced with\n# its corresponding new string, in document order. Because one string\n# (\"87\u00d761=\") occurs twice with two different replacements, each call only\n# replaces a single (the next, left-most still-unreplaced) occurrence, so\n# repeating the same \"old\" text with a different \"new\" text on the next line\n# correctly targets the next hit.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstOccurrence($findText, $replText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replText\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $found = $find.Execute($findText, $false, $true, $false, $false, $false, $true, 0, $false, $replText, 1)\n  if (-not $found) {\n    Write-Output \"Could not find text to replace: $findText\"\n  }\n}\n\nReplace-FirstOccurrence \"2024-08-22 Thursday\" \"2024-08-23 Friday\"\nReplace-FirstOccurrence \"38\u00d769=\" \"50\u00d722=\"\nReplace-FirstOccurrence \"65\u00d767=\" \"68\u00d793=\"\nReplace-FirstOccurrence \"22\u00d798=\" \"38\u00d733=\"\nReplace-FirstOccurrence \"45\u00d713=\" \"59\u00d728=\"\nReplace-FirstOccurrence \"38\u00d730=\" \"52\u00d770=\"\nReplace-FirstOccurrence \"48\u00d753=\" \"91\u00d742=\"\nReplace-FirstOccurrence \"23\u00d766=\" \"79\u00d766=\"\nReplace-FirstOccurrence \"80\u00d740=\" \"42\u00d798=\"\nReplace-FirstOccurrence \"27\u00d770=\" \"35\u00d725=\"\nReplace-FirstOccurrence \"35\u00d723=\" \"79\u00d766=\"\nReplace-FirstOccurrence \"54\u00d753=\" \"56\u00d720=\"\nReplace-FirstOccurrence \"87\u00d761=\" \"48\u00d743=\"\nReplace-FirstOccurrence \"64\u00d798=\" \"29\u00d758=\"\nReplace-FirstOccurrence \"87\u00d761=\" \"19\u00d727=\"\nReplace-FirstOccurrence \"35\u00d791=\" \"42\u00d776=\"\nReplace-FirstOccurrence \"72\u00d727=\" \"14\u00d788=\"\nReplace-FirstOccurrence \"95\u00d792=\" \"98\u00d777=\"\nReplace-FirstOccurrence \"12\u00d725=\" \"70\u00d776=\"\nReplace-FirstOccurrence \"43\u00d711=\" \"84\u00d724=\"\nReplace-FirstOccurrence \"20\u00d738=\" \"11\u00d794=\"\nReplace-FirstOccurrence \"19\u00d769=\" \"51\u00d759=\"\nReplace-FirstOccurrence \"23\u00d737=\" \"83\u00d751=\"\nReplace-FirstOccurrence \"52\u00d793=\" \"29\u00d748=\"\nReplace-FirstOccurrence \"42\u00d772=\" \"97\u00d746=\"\nReplace-FirstOccurrence \"51\u00d742=\" \"61\u00d736=\"\n"}
